$d = $word.ActiveDocument

# 1) "json" -> "excel" (also removes the spellcheck proofErr wrapper around the word)
$r = $d.Content
$r.Find.Execute("json", $false, $false, $false, $false, $false, $true, 1, $false, "excel", 2) | Out-Null
